$d = $word.ActiveDocument

# 1. Update "21 years" -> "15+ years" in professional summary
$d.Content.Find.Execute(
    "21 years of expertise", $true, $false, $false, $false, $false,
    $true, 1, $false, "15+ years of expertise", 2
) | Out-Null

# 2. Update FLEEM bullet (Progressive Change Campaign Committee section)
$d.Content.Find.Execute(
    "Conceived, architected, and engineered FLEEM web application using Twilio API for thousands of simultaneous phone calls",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Conceived, architected, and engineered FLEEM web application using Twilio API handling tens of thousands of calls using emulated predictive dialer for regulated political surveys",
    2
) | Out-Null

# 3. Salsa Labs section - replace first 4 bullets, then append 2 new bullets
$d.Content.Find.Execute(
    "Developed software solutions for political campaigns and advocacy groups",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Maintained and extended comprehensive geospatial analysis and reporting tools for Java-based CRM system used by tens of thousands of users simultaneously",
    2
) | Out-Null
$d.Content.Find.Execute(
    "Built web applications for voter engagement and campaign management",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Developed custom tile server for Web Map Service (WMS) integration using GeoTools and OpenLayers",
    2
) | Out-Null
$d.Content.Find.Execute(
    "Integrated third-party APIs and data sources for campaign tools",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Built advanced geospatial analysis capabilities using Java, JavaScript, MySQL, and TileMill",
    2
) | Out-Null
$d.Content.Find.Execute(
    "Collaborated with political strategists to translate requirements into technical solutions",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Integrated mapping and visualization tools for political campaign data analysis interfacing with Government and Activism APIs",
    2
) | Out-Null

# 4. Praxis Project section - replace first 4 bullets
$d.Content.Find.Execute(
    "Integrated technology solutions within organizational frameworks for social justice organizations",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Led technology operations for multi-million dollar organization while assisting in search for full-time CTO",
    2
) | Out-Null
$d.Content.Find.Execute(
    "Developed data management systems for community organizing efforts",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Directed all technology decisions and practices for massive multinational non-governmental organization",
    2
) | Out-Null
$d.Content.Find.Execute(
    "Provided technical training and support to nonprofit staff",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Developed comprehensive frameworks for internal and external technology audits",
    2
) | Out-Null
$d.Content.Find.Execute(
    "Built custom applications for community engagement and advocacy",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Led training initiatives for beneficiaries on spatial and Census data analysis for public health research",
    2
) | Out-Null

# --- Now handle insertions of brand-new bullet paragraphs. ---
# Process from the bottom of the document upward so paragraph indices
# for higher-up insertion points remain valid.

function Find-ParaIndexByText($doc, $substring) {
    $idx = 0
    $found = 0
    foreach ($para in $doc.Paragraphs) {
        $idx = $idx + 1
        if ($para.Range.Text.Contains($substring)) {
            $found = $idx
        }
    }
    return $found
}

# 6. Feldman Group: insert new bullet after "Enhanced value of research deliverables..."
$idx = Find-ParaIndexByText $d "Enhanced value of research deliverables through advanced analytical techniques using SPSS, OSCAR, PHP, and MySQL"
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter() | Out-Null
$newIdx = $idx + 1
$d.Paragraphs.Item($newIdx).Range.Text = "• Trained staff on PHP/MySQL for data analysis and reporting systems"

# 5. Lake Research Partners: insert new bullet after "Developed innovative approaches to visualizing..."
$idx = Find-ParaIndexByText $d "Developed innovative approaches to visualizing demographic and market data for enhanced client understanding"
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter() | Out-Null
$newIdx = $idx + 1
$d.Paragraphs.Item($newIdx).Range.Text = "• Trained staff on building Python tooling for report generation and analysis"

# 4b. Praxis Project: append 3 new bullets after the (already-replaced) 4th bullet
$idx = Find-ParaIndexByText $d "Led training initiatives for beneficiaries on spatial and Census data analysis for public health research"
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs.Item($idx + 1).Range.Text = "• Conducted training programs for NGO staff in web development using Drupal, PHP, and MySQL"
$d.Paragraphs.Item($idx + 1).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs.Item($idx + 2).Range.Text = "• Managed technology infrastructure supporting community health initiatives across multiple countries"
$d.Paragraphs.Item($idx + 2).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs.Item($idx + 3).Range.Text = "• Architected and developed 25 Drupal sites to integrate with membership databases, activism CRMs and government agencies, under guidelines from Kellogg Foundation and Robert Wood Johnson Foundation"

# 3b. Salsa Labs: append 2 new bullets after the (already-replaced) 4th bullet
$idx = Find-ParaIndexByText $d "Integrated mapping and visualization tools for political campaign data analysis interfacing with Government and Activism APIs"
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs.Item($idx + 1).Range.Text = "• Collaborated with political strategists to translate geospatial requirements into technical solutions"
$d.Paragraphs.Item($idx + 1).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs.Item($idx + 2).Range.Text = "• Handled billions of records with millions of columns in high-performance CRM system"

Write-Output "Edit complete"
